$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A32").Value = "Globo"
$ws.Range("B32").Value = "RJ TV 2"
$ws.Range("C32").Value = "Transporte"
$ws.Range("D32").Value = "2025-04-01T19:18"
$ws.Range("E32").Value = "Negativo"
$ws.Range("F32").Value = "Mais um dia sem vans. Moradores têm problemas para voltar para casa depois de paralisação do setor C. Em 15 dias, é a segunda paralisação. Entrevista com passageiros reclamando. Problema começou segunda-feira. Vans do setor C (que atende região Norte) pararam 100%. Mês passado, foi 50% de paralisação. Alguns passageiros estão tendo que pagar R$ 20 para conseguir voltar para casa. Entrevista com persmissionário, Jefferson Henrique. Relatórios estariam com inconsistência. Eles receberam menos que o previsto. *com nota* do IMTT"
